$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Re-ran cross-validation with a new k; update the reported metrics for the
# random_forest, lsboost and neural_network rows (old_model / row 5 and the
# header row / row labels are unchanged).

$ws.Range("B2").Value = 50.912618681310462
$ws.Range("C2").Value = 0.43282709726054003
$ws.Range("D2").Value = 38.469694995520697
$ws.Range("E2").Value = 0.36903534952664963
$ws.Range("F2").Value = 0.60748279772076641
$ws.Range("G2").Value = 0.60876457539570772
$ws.Range("H2").Value = 0.63096465047335037
$ws.Range("I2").Value = 0.8016636574855327

$ws.Range("B3").Value = 49.417167662041138
$ws.Range("C3").Value = 0.42011371223869937
$ws.Range("D3").Value = 36.811477558999442
$ws.Range("E3").Value = 0.34767446765454485
$ws.Range("F3").Value = 0.58963926909131892
$ws.Range("G3").Value = 0.58252407534040018
$ws.Range("H3").Value = 0.65232553234545509
$ws.Range("I3").Value = 0.81400553741348491

$ws.Range("B4").Value = 54.401817069252665
$ws.Range("C4").Value = 0.46249006980321078
$ws.Range("D4").Value = 39.202741455402894
$ws.Range("E4").Value = 0.42135088396048764
$ws.Range("F4").Value = 0.64911546273408682
$ws.Range("G4").Value = 0.6203646860008859
$ws.Range("H4").Value = 0.57864911603951241
$ws.Range("I4").Value = 0.76941770635881213

